$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

# ---- Title ----
Replace-Text "Einstein's Relativity - Redefining Space and Time" "The Mathematical Realm: Unraveling Patterns and Structures"

# ---- Author name ----
Replace-Text "Elijah Stone" "Clara Hudson"

# ---- Email parts ----
Replace-Text "elijah" "clara"
Replace-Text "stone@xyz" "hudson@xyz"

# ---- Body paragraph 1 (three sentences) ----
Replace-Text "In the annals of scientific history, Albert Einstein's revolutionary theory of relativity stands as a towering testament to the transformative power of human intellect" "Delving into mathematics is akin to embarking on an enthralling voyage of discovery, where the interplay of numbers, patterns, and structures unfolds before our very eyes"

Replace-Text " This profound intellectual odyssey, spanning the early 20th century, unveiled a paradigm-shifting understanding of space, time, gravity, and the underlying fabric of the universe" " It is a field that captivates the mind with its precise reasoning and abstract beauty, nurturing analytical thinking and problem-solving skills that are indispensable in our ever-evolving world"

Replace-Text " Einstein's groundbreaking work challenged long-held classical notions, upending our fundamental perceptions of reality and ushering in a new era of scientific discovery" " From the ancient civilizations that deciphered numerical systems to the modern-day advancements in computer science, mathematics has left an indomitable mark on human progress and continues to shape our understanding of the universe"

# ---- Body paragraph 2 (four sentences) ----
Replace-Text "In his seminal 1905 paper, Einstein introduced the concept of special relativity, delving into the intricate relationship between space and time" "Mathematics is not merely a collection of abstract concepts; it finds practical applications in numerous disciplines, enriching our lives in countless ways"

Replace-Text " His groundbreaking postulates shattered the traditional view of absolute time and distance, revealing that these concepts are relative to the observer's frame of reference" " Engineers utilize mathematical principles to design and construct robust structures, while economists leverage mathematical models to analyze market trends and predict economic behavior"

Replace-Text " This radical departure from classical physics led to the iconic equation E=mc2, which elegantly encapsulates the equivalence between mass and energy" " Mathematicians themselves engage in groundbreaking research, pushing the boundaries of knowledge and expanding our comprehension of the cosmos"

Replace-Text " Einstein's special theory of relativity transformed our understanding of the universe's fundamental building blocks" " Its universality transcends cultural and linguistic barriers, connecting individuals from all corners of the globe in a shared pursuit of understanding"

# ---- Body paragraph 3 (three sentences) ----
Replace-Text "Einstein's intellectual journey culminated in the formulation of general relativity, a profound theory that elucidated the nature of gravity and redefined our perception of the cosmos" "The beauty of mathematics lies in its elegance and simplicity, where complex phenomena can be distilled into concise equations or formulas"

Replace-Text " General relativity conceptualizes gravity not as a force but as a curvature of spacetime caused by the presence of mass and energy" " It is a language that describes the fundamental workings of the universe, capable of expressing the laws of physics, the intricacies of biology, and the patterns of human behavior"

Replace-Text " This elegant framework revolutionized our understanding of celestial phenomena, explaining the intricate motions of planets and galaxies and opening up new avenues for exploring the vast expanse of the universe" " By delving into this realm of numbers and relationships, we uncover hidden truths and gain a deeper appreciation for the order and harmony that underpin our existence"

# ---- Summary paragraph ----
# The Summary paragraph undergoes a structural change: the tail sentences are
# replaced, one run is split in two (with a lastRenderedPageBreak marker on
# the second part), and a trailing sentence is added. We rebuild the whole
# paragraph with InsertXML so the exact run layout (including the page-break
# marker) is preserved, instead of relying on Find/Replace (which normalizes
# same-formatted adjacent runs into a single run).
$summaryPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$summaryRange = $summaryPara.Range

$summaryXml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>Mathematics, a captivating field of study, invites us to unravel patterns, explore structures, and uncover hidden truths</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Its practical applications span a multitude of disciplines, while its aesthetic elegance captivates the mind</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> As we delve deeper into the mathematical realm, we </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>cultivate analytical thinking, problem-solving skills, and a profound appreciation for the order and harmony that govern our universe</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Mathematics transcends cultural and linguistic boundaries, serving as a universal language that unites individuals in a shared pursuit of understanding and progress</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$summaryRange.InsertXML($summaryXml)

Write-Host "Done"
